$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "Discount" line inside the first merged group (row 5): the
# discount amount in C5, and an italic+underlined "Discount" label in B5.
$ws.Range("C5").Value = -11.41

$ws.Range("B5").Value = "Discount"
$ws.Range("B5").Font.Italic = $true
$ws.Range("B5").Font.Underline = $true

# Move the active selection like the author's saved session.
$ws.Range("L4").Select()
